$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Formula = "=""67.488.55"""
$ws.Cells.Item(2, 4).Copy()
$ws.Cells.Item(2, 4).PasteSpecial(-4163)
$ws.Cells.Item(2, 5).Value = "  +3.32%  "
$ws.Cells.Item(3, 4).Formula = "=""3.279.71"""
$ws.Cells.Item(3, 4).Copy()
$ws.Cells.Item(3, 4).PasteSpecial(-4163)
$ws.Cells.Item(3, 5).Value = "  -0.13%  "
$ws.Cells.Item(4, 5).Value = "  -0.01%  "
$ws.Cells.Item(5, 4).Formula = "=""573.93"""
$ws.Cells.Item(5, 4).Copy()
$ws.Cells.Item(5, 4).PasteSpecial(-4163)
$ws.Cells.Item(5, 5).Value = "  -0.51%  "
$ws.Cells.Item(6, 4).Formula = "=""176.32"""
$ws.Cells.Item(6, 4).Copy()
$ws.Cells.Item(6, 4).PasteSpecial(-4163)
$ws.Cells.Item(6, 5).Value = "  -2.88%  "
$ws.Cells.Item(7, 5).Value = "  +0.15%  "
$ws.Cells.Item(8, 4).Formula = "=""0.583"""
$ws.Cells.Item(8, 4).Copy()
$ws.Cells.Item(8, 4).PasteSpecial(-4163)
$ws.Cells.Item(8, 5).Value = "  +2.93%  "
$ws.Cells.Item(9, 4).Formula = "=""3.270.74"""
$ws.Cells.Item(9, 4).Copy()
$ws.Cells.Item(9, 4).PasteSpecial(-4163)
$ws.Cells.Item(9, 5).Value = "  -0.19%  "
$ws.Cells.Item(10, 5).Value = "  -0.10%  "
$ws.Cells.Item(11, 4).Formula = "=""0.572"""
$ws.Cells.Item(11, 4).Copy()
$ws.Cells.Item(11, 4).PasteSpecial(-4163)
$ws.Cells.Item(11, 5).Value = "  +1.23%  "
$ws.Cells.Item(12, 4).Formula = "=""45.36"""
$ws.Cells.Item(12, 4).Copy()
$ws.Cells.Item(12, 4).PasteSpecial(-4163)
$ws.Cells.Item(12, 5).Value = "  -1.64%  "
$ws.Cells.Item(13, 4).Formula = "=""0.0000268"""
$ws.Cells.Item(13, 4).Copy()
$ws.Cells.Item(13, 4).PasteSpecial(-4163)
$ws.Cells.Item(13, 5).Value = "  +2.54%  "
$ws.Cells.Item(14, 4).Formula = "=""701.94"""
$ws.Cells.Item(14, 4).Copy()
$ws.Cells.Item(14, 4).PasteSpecial(-4163)
$ws.Cells.Item(14, 5).Value = "  +12.66%  "
$ws.Cells.Item(15, 4).Formula = "=""3.817.41"""
$ws.Cells.Item(15, 4).Copy()
$ws.Cells.Item(15, 4).PasteSpecial(-4163)
$ws.Cells.Item(15, 5).Value = "  +0.48%  "
$ws.Cells.Item(16, 4).Formula = "=""8.32"""
$ws.Cells.Item(16, 4).Copy()
$ws.Cells.Item(16, 4).PasteSpecial(-4163)
$ws.Cells.Item(16, 5).Value = "  -0.06%  "
$ws.Cells.Item(17, 4).Formula = "=""67.576.40"""
$ws.Cells.Item(17, 4).Copy()
$ws.Cells.Item(17, 4).PasteSpecial(-4163)
$ws.Cells.Item(17, 5).Value = "  +3.18%  "
$ws.Cells.Item(18, 5).Value = "  +1.44%  "
$ws.Cells.Item(19, 4).Formula = "=""3.296.36"""
$ws.Cells.Item(19, 4).Copy()
$ws.Cells.Item(19, 4).PasteSpecial(-4163)
$ws.Cells.Item(19, 5).Value = "  +0.48%  "
$ws.Cells.Item(20, 4).Formula = "=""17.34"""
$ws.Cells.Item(20, 4).Copy()
$ws.Cells.Item(20, 4).PasteSpecial(-4163)
$ws.Cells.Item(20, 5).Value = "  -1.23%  "
$ws.Cells.Item(21, 4).Formula = "=""10.73"""
$ws.Cells.Item(21, 4).Copy()
$ws.Cells.Item(21, 4).PasteSpecial(-4163)
$ws.Cells.Item(21, 5).Value = "  -0.94%  "
$ws.Cells.Item(22, 4).Formula = "=""0.888"""
$ws.Cells.Item(22, 4).Copy()
$ws.Cells.Item(22, 4).PasteSpecial(-4163)
$ws.Cells.Item(22, 5).Value = "  +0.80%  "
$ws.Cells.Item(23, 5).Value = "  -6.66%  "
$ws.Cells.Item(24, 4).Formula = "=""5.13"""
$ws.Cells.Item(24, 4).Copy()
$ws.Cells.Item(24, 4).PasteSpecial(-4163)
$ws.Cells.Item(24, 5).Value = "  +3.98%  "
$ws.Cells.Item(25, 4).Formula = "=""98.92"""
$ws.Cells.Item(25, 4).Copy()
$ws.Cells.Item(25, 4).PasteSpecial(-4163)
$ws.Cells.Item(25, 5).Value = "  -0.27%  "
$ws.Cells.Item(26, 4).Formula = "=""3.90"""
$ws.Cells.Item(26, 4).Copy()
$ws.Cells.Item(26, 4).PasteSpecial(-4163)
$ws.Cells.Item(26, 5).Value = "  -0.62%  "
$ws.Cells.Item(27, 4).Formula = "=""2.71"""
$ws.Cells.Item(27, 4).Copy()
$ws.Cells.Item(27, 4).PasteSpecial(-4163)
$ws.Cells.Item(27, 5).Value = "  +0.27%  "
$ws.Cells.Item(28, 4).Formula = "=""9.27"""
$ws.Cells.Item(28, 4).Copy()
$ws.Cells.Item(28, 4).PasteSpecial(-4163)
$ws.Cells.Item(28, 5).Value = "  -0.38%  "
$ws.Cells.Item(29, 4).Formula = "=""32.85"""
$ws.Cells.Item(29, 4).Copy()
$ws.Cells.Item(29, 4).PasteSpecial(-4163)
$ws.Cells.Item(29, 5).Value = "  +7.82%  "
$ws.Cells.Item(30, 4).Formula = "=""8.41"""
$ws.Cells.Item(30, 4).Copy()
$ws.Cells.Item(30, 4).PasteSpecial(-4163)
$ws.Cells.Item(30, 5).Value = "  +1.35%  "
$ws.Cells.Item(31, 4).Formula = "=""6.63"""
$ws.Cells.Item(31, 4).Copy()
$ws.Cells.Item(31, 4).PasteSpecial(-4163)
$ws.Cells.Item(31, 5).Value = "  +3.04%  "
$ws.Cells.Item(32, 4).Formula = "=""577.94"""
$ws.Cells.Item(32, 4).Copy()
$ws.Cells.Item(32, 4).PasteSpecial(-4163)
$ws.Cells.Item(32, 5).Value = "  +5.32%  "
$ws.Cells.Item(33, 4).Formula = "=""3.884.62"""
$ws.Cells.Item(33, 4).Copy()
$ws.Cells.Item(33, 4).PasteSpecial(-4163)
$ws.Cells.Item(33, 5).Value = "  +2.22%  "
$ws.Cells.Item(34, 4).Formula = "=""10.78"""
$ws.Cells.Item(34, 4).Copy()
$ws.Cells.Item(34, 4).PasteSpecial(-4163)
$ws.Cells.Item(34, 5).Value = "  +0.19%  "
$ws.Cells.Item(35, 5).Value = "  +0.42%  "
$ws.Cells.Item(36, 4).Formula = "=""0.999"""
$ws.Cells.Item(36, 4).Copy()
$ws.Cells.Item(36, 4).PasteSpecial(-4163)
$ws.Cells.Item(36, 5).Value = "  -0.04%  "
$ws.Cells.Item(37, 2).Value = "OKB"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(37, 4).Formula = "=""55.30"""
$ws.Cells.Item(37, 4).Copy()
$ws.Cells.Item(37, 4).PasteSpecial(-4163)
$ws.Cells.Item(37, 5).Value = "  -0.71%  "
$ws.Cells.Item(38, 2).Value = "dogwifhat"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Cells.Item(38, 4).Formula = "=""3.31"""
$ws.Cells.Item(38, 4).Copy()
$ws.Cells.Item(38, 4).PasteSpecial(-4163)
$ws.Cells.Item(38, 5).Value = "  -8.65%  "
$ws.Cells.Item(39, 5).Value = "  +2.58%  "
$ws.Cells.Item(40, 2).Value = "Fetch.AI"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(40, 4).Formula = "=""2.60"""
$ws.Cells.Item(40, 4).Copy()
$ws.Cells.Item(40, 4).PasteSpecial(-4163)
$ws.Cells.Item(40, 5).Value = "  +1.38%  "
$ws.Cells.Item(41, 2).Value = "Stacks"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(41, 4).Formula = "=""3.11"""
$ws.Cells.Item(41, 4).Copy()
$ws.Cells.Item(41, 4).PasteSpecial(-4163)
$ws.Cells.Item(41, 5).Value = "  -0.01%  "
$ws.Cells.Item(42, 4).Formula = "=""31.79"""
$ws.Cells.Item(42, 4).Copy()
$ws.Cells.Item(42, 4).PasteSpecial(-4163)
$ws.Cells.Item(42, 5).Value = "  -1.21%  "
$ws.Cells.Item(43, 2).Value = "ApeXProtocol"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Cells.Item(43, 4).Formula = "=""3.33"""
$ws.Cells.Item(43, 4).Copy()
$ws.Cells.Item(43, 4).PasteSpecial(-4163)
$ws.Cells.Item(43, 5).Value = "  -2.39%  "
$ws.Cells.Item(44, 2).Value = "PEPE"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Cells.Item(44, 4).Formula = "=""0.0₃0672"""
$ws.Cells.Item(44, 4).Copy()
$ws.Cells.Item(44, 4).PasteSpecial(-4163)
$ws.Cells.Item(44, 5).Value = "  +0.17%  "
$ws.Cells.Item(45, 4).Formula = "=""0.328"""
$ws.Cells.Item(45, 4).Copy()
$ws.Cells.Item(45, 4).PasteSpecial(-4163)
$ws.Cells.Item(45, 5).Value = "  +0.02%  "
$ws.Cells.Item(46, 4).Formula = "=""0.0406"""
$ws.Cells.Item(46, 4).Copy()
$ws.Cells.Item(46, 4).PasteSpecial(-4163)
$ws.Cells.Item(46, 5).Value = "  +1.00%  "
$ws.Cells.Item(47, 5).Value = "  +1.86%  "
$ws.Cells.Item(48, 4).Formula = "=""1.39"""
$ws.Cells.Item(48, 4).Copy()
$ws.Cells.Item(48, 4).PasteSpecial(-4163)
$ws.Cells.Item(48, 5).Value = "  +11.13%  "
$ws.Cells.Item(49, 5).Value = "  +0.33%  "
$ws.Cells.Item(50, 4).Formula = "=""2.51"""
$ws.Cells.Item(50, 4).Copy()
$ws.Cells.Item(50, 4).PasteSpecial(-4163)
$ws.Cells.Item(50, 5).Value = "  +1.00%  "
$ws.Cells.Item(51, 4).Formula = "=""127.91"""
$ws.Cells.Item(51, 4).Copy()
$ws.Cells.Item(51, 4).PasteSpecial(-4163)
$ws.Cells.Item(51, 5).Value = "  -0.10%  "
$excel.CutCopyMode = 0
